$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 139.3243227574293
$ws.Range("C3").Value = 10.58980444494598
$ws.Range("C4").Value = 6.854514699577449
$ws.Range("C5").Value = 16.1526472176978
$ws.Range("C6").Value = 40.79285717030207
$ws.Range("C7").Value = 12.59598450584848
$ws.Range("C8").Value = 9.014710958269855
$ws.Range("C9").Value = 27.79754200729729
$ws.Range("C10").Value = 44.83459151539689
$ws.Range("C11").Value = 10.52867430458607
$ws.Range("C12").Value = 4.015535713512015
$ws.Range("C13").Value = 7.314975497093655
$ws.Range("C14").Value = 1.868835719574394
$ws.Range("C15").Value = 1.808499477141236
$ws.Range("C16").Value = 20.90571410516249
$ws.Range("C17").Value = 19.02020652912631
$ws.Range("C18").Value = 10.19920666498396
$ws.Range("C19").Value = 1.085258465870092
$ws.Range("C20").Value = 30.41264177801864
$ws.Range("C21").Value = 78.40297755225505
$ws.Range("C22").Value = 13.90512218706266
$ws.Range("C23").Value = 0.22149752156383
$ws.Range("C24").Value = 2.527770998778619
$ws.Range("C25").Value = 26.90917022726171
$ws.Range("C26").Value = 7.422151717205185
$ws.Range("C27").Value = 0.5795454865290176
$ws.Range("C28").Value = 10.66601864591418
$ws.Range("C29").Value = 25.06176975170963
$ws.Range("C30").Value = 10.10314501584696
$ws.Range("C31").Value = 13.06994156811947
$ws.Range("C32").Value = 3.64716707549905
$ws.Range("C33").Value = 1.848988271405592
$ws.Range("C34").Value = 5.073007751945783
$ws.Range("C35").Value = 2.480137123173494
$ws.Range("C36").Value = 90.23126276293428
$ws.Range("C37").Value = 8.446280042715365
$ws.Range("C38").Value = 25.85011039297444
$ws.Range("C39").Value = 5.311971027898159
$ws.Range("C40").Value = 3.412173289180435
$ws.Range("C41").Value = 12.50944963183251
$ws.Range("C42").Value = 0.8891656779623284
$ws.Range("C43").Value = 5.836737557481284
$ws.Range("C44").Value = 230.477915190917
